# Rename the embedded logo pictures so the wp:docPr/@name (and, where the
# object model allows it, pic:cNvPr/@name) attributes swap their filenames:
#   - Footer 1 / Footer 2 : the Pearson logo  "image2.png" -> "image1.png"
#   - Header 2            : the BTec logo     "image1.jpg" -> "image2.jpg"
#
# InlineShape objects living in a header/footer story need to be re-fetched
# through their own (narrow) Range before the rename is applied - going
# through the wider Header/Footer Range directly can leave the shape's
# owning paragraph unresolved, so we always hop through shape.Range first.

$d   = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-FirstInlineShape($range, [string]$newName) {
    $shape = $range.InlineShapes.Item(1)
    # Re-acquire the shape via its own atomic Range so the rename resolves
    # correctly even inside a multi-paragraph header/footer story.
    $shape = $shape.Range.InlineShapes.Item(1)
    $shape.Name = $newName
}

# --- Footer 1: Pearson logo (id=1) ---
$footer1 = $sec.Footers.Item(1)
Rename-FirstInlineShape $footer1.Range "image1.png"

# --- Footer 2: Pearson logo (id=2) ---
$footer2 = $sec.Footers.Item(2)
Rename-FirstInlineShape $footer2.Range "image1.png"

# --- Header 2: BTec logo (id=3) ---
$header2 = $sec.Headers.Item(2)
Rename-FirstInlineShape $header2.Range "image2.jpg"
